# initial malware check implement
#
# Inserts a new "check" / "CHECK" localization row right after the
# existing "close" row (new row 10), shifting all subsequent rows down
# by one, and appends six new rows at the bottom of the table for the
# malware-check feature's localization strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes old rows 10-46 down to 11-47)
# and populate it with the new "check" key/value pair.
$ws.Rows(10).Insert()
$ws.Range("A10").Value = "check"
$ws.Range("B10").Value = "CHECK"

# Append the new malware-check localization rows at the end of the table
# (now rows 48-53).
$ws.Range("A48").Value = "malware_check_title"
$ws.Range("B48").Value = "Malware Identifier"

$ws.Range("A49").Value = "malware_not_found"
$ws.Range("B49").Value = "No malware found."

$ws.Range("A50").Value = "malware_trojan_title"
$ws.Range("B50").Value = "trojan"

$ws.Range("A51").Value = "malware_trojan_detail"
$ws.Range("B51").Value = "trojan detail"

$ws.Range("A52").Value = "malware_rat_title"
$ws.Range("B52").Value = "rat"

$ws.Range("A53").Value = "malware_rat_detail"
$ws.Range("B53").Value = "rat detail"

# Reflect the final cursor / selection position left behind by the edit.
$ws.Range("B44").Select()
